$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.292.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -7.39%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.465.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.63%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -1.11%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'380.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -8.69%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'122.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -5.79%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'3.562.24"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.71%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.572"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -12.41%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.996"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.34%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.641"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -17.52%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.141"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -20.09%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.0000290"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -14.06%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'38.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -10.67%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'3.952.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -5.10%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'9.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -8.94%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.136"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.10%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.402.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -5.17%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.36%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'18.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -11.61%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "'61.488.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -8.57%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("B21").Value = "Polygon"
$ws.Range("C21").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D21").Value = "'0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -12.62%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'377.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -16.71%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'13.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.01%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'78.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -12.19%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -14.17%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'5.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +5.20%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'32.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -7.05%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'2.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -13.80%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'8.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -13.97%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'11.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -6.86%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'2.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -7.36%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.107"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -9.31%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'6.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -14.41%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'55.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.32%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").Value = "'0.996"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.41%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.142"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -12.96%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "'35.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -12.28%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.0421"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -14.83%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -0.52%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.129"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -12.15%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").Value = "'2.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +12.91%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'139.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -6.46%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'25.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +19.14%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'2.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +12.64%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'2.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -11.69%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "'1.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.12%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "PEPE"
$ws.Range("C47").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D47").Value = "'0.0₃0576"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -23.57%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").Value = "'2.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -9.40%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'3.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -9.43%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -17.63%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.262"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -17.92%  "
$ws.Range("E51").Style = "Normal"
